$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: replace "FREDDY CASTILLO FALCON" with "KIRA PLANCHART" and update the C.I. number.
# The new name also needs the Arial 10pt font (no special alignment) - the same
# font already used elsewhere in the sheet (e.g. B12), but WITHOUT the centered
# alignment that those cells carry. Build that exact format on a scratch cell
# first (so the alignment stays "general"/untouched) and then copy *only the
# formatting* onto A11 - this reuses the existing font instead of minting a
# brand-new one.
$helper = $ws.Range("Z1")
$helper.Value = "helper"
$helper.Font.Size = 10
$helper.Font.Name = "Arial"

$helper.Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$helper.Delete(-4159)

$ws.Range("A11").Value = "KIRA PLANCHART"
$ws.Range("B11").Value = 16954727

# Row 6: the "CARGO" for this employee changed.
$ws.Range("C6").Value = "SUPERVISOR ALEDAÑOS"

# Update the last selected cell shown when the workbook is reopened.
$ws.Range("A16").Select()
